$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells in columns D and E may look numeric (e.g. "1.00", "0.538", percentages).
# Force them to be treated/stored as literal text, matching the source data,
# then reset the style so no extra formatting is left behind on the cell.
function Set-TextValue($ws, $ref, $val) {
    $rng = $ws.Range($ref)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

Set-TextValue $ws 'D2' '47.347.48'
Set-TextValue $ws 'E2' '  +2.60%  '
Set-TextValue $ws 'D3' '2.505.94'
Set-TextValue $ws 'E3' '  +2.36%  '
Set-TextValue $ws 'E4' '  +0.05%  '
Set-TextValue $ws 'D5' '324.12'
Set-TextValue $ws 'E5' '  +0.88%  '
Set-TextValue $ws 'D6' '109.49'
Set-TextValue $ws 'E6' '  +4.70%  '
Set-TextValue $ws 'E7' '  +1.55%  '
Set-TextValue $ws 'E8' '  -0.03%  '
Set-TextValue $ws 'D9' '0.538'
Set-TextValue $ws 'E9' '  +0.16%  '
Set-TextValue $ws 'E10' '  +9.54%  '
Set-TextValue $ws 'D11' '0.0814'
Set-TextValue $ws 'E11' '  +1.15%  '
Set-TextValue $ws 'E12' '  +0.86%  '
Set-TextValue $ws 'D13' '18.46'
Set-TextValue $ws 'E13' '  +1.18%  '
Set-TextValue $ws 'E14' '  +2.21%  '
Set-TextValue $ws 'D15' '2.897.48'
Set-TextValue $ws 'E15' '  +2.24%  '
Set-TextValue $ws 'D16' '2.499.56'
Set-TextValue $ws 'E16' '  +2.45%  '
Set-TextValue $ws 'E17' '  +1.81%  '
Set-TextValue $ws 'D18' '47.276.61'
Set-TextValue $ws 'E18' '  +2.67%  '
Set-TextValue $ws 'D19' '12.90'
Set-TextValue $ws 'E19' '  +2.50%  '
Set-TextValue $ws 'D20' '6.66'
Set-TextValue $ws 'E20' '  +4.12%  '
Set-TextValue $ws 'E21' '  +1.28%  '
Set-TextValue $ws 'D22' '2.72'
Set-TextValue $ws 'E22' '  +14.66%  '
Set-TextValue $ws 'E23' '  -0.59%  '
Set-TextValue $ws 'D24' '248.21'
Set-TextValue $ws 'E24' '  +0.60%  '
Set-TextValue $ws 'D25' '2.61'
Set-TextValue $ws 'E25' '  +3.65%  '
Set-TextValue $ws 'D26' '26.07'
Set-TextValue $ws 'E26' '  +0.74%  '
Set-TextValue $ws 'D27' '1.00'
Set-TextValue $ws 'E27' '  -0.04%  '
Set-TextValue $ws 'D28' '2.29'
Set-TextValue $ws 'E28' '  +0.71%  '
Set-TextValue $ws 'D29' '10.03'
Set-TextValue $ws 'E29' '  +3.89%  '
Set-TextValue $ws 'D30' '35.77'
Set-TextValue $ws 'E30' '  +5.73%  '
Set-TextValue $ws 'E31' '  +8.40%  '
Set-TextValue $ws 'D32' '49.87'
Set-TextValue $ws 'E32' '  +0.81%  '
Set-TextValue $ws 'D33' '19.97'
Set-TextValue $ws 'E33' '  +0.68%  '
Set-TextValue $ws 'D34' '5.45'
Set-TextValue $ws 'E34' '  +2.16%  '
Set-TextValue $ws 'E35' '  +4.03%  '
Set-TextValue $ws 'E36' '  +0.24%  '
Set-TextValue $ws 'E37' '  +5.42%  '
Set-TextValue $ws 'D38' '4.72'
Set-TextValue $ws 'E38' '  +4.60%  '
Set-TextValue $ws 'E39' '  +1.85%  '
Set-TextValue $ws 'E40' '  +1.25%  '
Set-TextValue $ws 'E41' '  +0.54%  '
Set-TextValue $ws 'D42' '120.68'
Set-TextValue $ws 'E42' '  -4.94%  '
Set-TextValue $ws 'D43' '21.21'
Set-TextValue $ws 'E43' '  +2.68%  '
Set-TextValue $ws 'E44' '  +2.36%  '
Set-TextValue $ws 'D45' '2.001.61'
Set-TextValue $ws 'E45' '  +1.94%  '
Set-TextValue $ws 'D46' '3.09'
Set-TextValue $ws 'E46' '  +4.34%  '
Set-TextValue $ws 'E47' '  -0.51%  '
$ws.Range('B48').Value = 'FraxShare'
$ws.Range('C48').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
Set-TextValue $ws 'D48' '9.10'
Set-TextValue $ws 'E48' '  +0.05%  '
$ws.Range('B49').Value = 'Stacks'
$ws.Range('C49').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
Set-TextValue $ws 'D49' '1.78'
Set-TextValue $ws 'E49' '  -4.19%  '
Set-TextValue $ws 'E50' '  +4.40%  '
Set-TextValue $ws 'D51' '56.97'
Set-TextValue $ws 'E51' '  +4.67%  '
